$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column)
$ws.Columns("N:N").Insert()

# Match the width Excel assigns the newly inserted column (copied from the
# column to its left rather than being re-measured with AutoFit)
$ws.Columns("N:N").ColumnWidth = 9.83

# Select cell R5 as the active cell on this sheet
$ws.Range("R5").Select()
